$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 2002
$ws.Range("B5").Value = "TRPGスーパーセッション大饗宴"
$ws.Range("C5").Value = "RPG Super Session Feast"
$ws.Range("D5").Value = "Enterbrain"
$ws.Range("E5").Value = "rpg_super_session_feast.jpg"
$ws.Range("F5").Value = "periodical"

$ws.Columns.Item(2).ColumnWidth = 33.833333333333336
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(5).ColumnWidth = 30.166666666666668

$ws.Range("A6").Select()
